$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 4: replace "pepito perez" record with "Panfilo" record
$ws.Cells.Item(4,1).Value = "Panfilo"
$ws.Cells.Item(4,2).Value = "webos"
$ws.Cells.Item(4,3).NumberFormat = "@"
$ws.Cells.Item(4,3).Value = "12345678"
$ws.Cells.Item(4,3).NumberFormat = "General"
$ws.Cells.Item(4,3).Style = "Normal"
$ws.Cells.Item(4,4).Value = "mailserio@mail.com"
$ws.Cells.Item(4,5).Value = 1
$ws.Cells.Item(4,6).Value = 3

# Row 7 (new): "manuel" record appended at the end
$ws.Cells.Item(7,1).Value = "manuel"
$ws.Cells.Item(7,2).Value = "manuel"
$ws.Cells.Item(7,3).NumberFormat = "@"
$ws.Cells.Item(7,3).Value = "12345678"
$ws.Cells.Item(7,3).NumberFormat = "General"
$ws.Cells.Item(7,3).Style = "Normal"
$ws.Cells.Item(7,4).Value = "hola@hola.com"
$ws.Cells.Item(7,5).Value = 1
$ws.Cells.Item(7,6).Value = 6

[void]$ws.Range("E5").Select()
